$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.716.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.438.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.12%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.438.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("E9").Value = "  -7.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "

$ws.Range("E11").Value = "  -4.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.030.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.776.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.440.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.541"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.98%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "

$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("E32").Value = "  -2.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.899.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0744"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.779"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0309"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "317.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.07%  "

$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.80%  "

$ws.Range("E51").Value = "  -3.59%  "
